# Fix Training Data Issue (#48)
# The "Date" column (BF) held the sheet's own filename-derived label
# ("5-30-2012-13") instead of an actual date string. NBA.com displayed
# the stats for this date one day off from the real date, so the value
# is corrected to the proper ISO-style date "2013-05-30" for every data
# row (rows 2-31; row 1 is the "Date" header and must stay untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$dateCol = 58   # column BF

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $dateCol)
    # Force plain text so "2013-05-30" is stored as a literal string
    # instead of being auto-converted into a date serial number.
    $cell.NumberFormat = "@"
    $cell.Value = "2013-05-30"
}
